# Consolidate runs that were split purely by spell-check <w:proofErr> markers
# (the visible text is unchanged; this merges adjacent runs with identical
# formatting back into a single run, matching the author's cleanup).
$d = $word.ActiveDocument

$mergeTexts = @(
    "Notizen von mkerk in rot",
    "Entwicklung einer App zur Passwortverwaltung",
    "Mobil (Deshalb eine App für’s Smartphone)",
    "Nicht nur Verwaltung sondern auch Nutzung soll erleichtert werden (Copy-Funktion)",
    "Stabilität der App soll gewährleistet sein",
    "Ein weiteres Ziel des Projekts war die Einarbeitung in die Entwicklung von Android Apps",
    "Wurde nicht umgesetzt, da die Daten zur Generierung von Passwörtern über die Laufzeit der App hinaus in eine zusätz",
    "In diesem Fall ist ein Brute-Force Angriff oder sogar ein Wörterbuchangriff durch raten des Masterpassworts möglich",
    "Die App ist für die Verwendung auf einem Smartphone optimiert",
    "Da Tablets weniger mobil sind wurde hier der Schwerpunkt gelegt",
    "Ggf. UML-Diagramme, Screenshots",
    "Es wurde ein inkrementelles Vorgehen zur Entwicklung der App gewählt",
    "So entstehen keine technischen Abhängigkeiten und der Blick auf die geforderten Anforderungen geht nicht verlorgen",
    "Anwendung eines inkrementellen Vorgehens da wir den Aufwand zur Entwicklung von Android Apps anfänglich nicht beurteilen konnten",
    "Regelmäßige Absprachen ähnlich den Planning Meetings im Vorgehensmodell Scrum",
    "Git als VCS hatte für unsere Gruppe hauptsächlich den Vorteil, dass es ein zentral verfügbares Repository gab in dem wir unsere Entwicklungsstände zusammenführen konnten",
    "Paralleles Arbeiten fand nur in sehr seltenen Fällen statt, wodurch fast keine Merge Konflikte auftraten",
    "Die Möglichkeiten von Git als DVCS wurden nicht in Anspruch genommen, da dies den anfänglichen Einstieg nur erschwert hätte (Stichwort: Commit)",
    "Aufbauend auf Git als VCS wurde als grafische Oberfläche SourceTree verwendet",
    "Weitesgehend intuitive Nutzung",
    "Leider wird kein Merge Tool mit ausgeliefert (Allerdings kann mit ein wenig Aufwand ein externes Merge Tool angebunden werden)",
    "Die anderen Funktionen der App wurden durch Tests nach der Umsetzung der einzelnen Funktionen getestet",
)

foreach ($txt in $mergeTexts) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $ok = $range.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, $txt, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $txt"
    }
}

# New bullet point added by the author: a third-level item right after the
# "Leider wird kein Merge Tool..." bullet and before the "Da die App..." bullet.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Merge Tool angebunden werden)*") {
        $anchor = $p
    }
}
if ($null -eq $anchor) {
    Write-Output "ANCHOR PARAGRAPH NOT FOUND"
} else {
    $anchor.Range.InsertParagraphAfter()
    $newPar = $anchor.Next()
    $newPar.Range.Text = "Die Nutzung von Git Flow zur automatisierten Verwaltung der Branches schlug fehl, da dies in Verbindung mit dem neuen VCS wahrscheinlich Verwirrung schaffte"
}

Write-Output "Done"
